# Auto-generated edit script applying the crypto price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'70.761.32"
$ws.Range("E2").Value = "  +2.77%  "
$ws.Range("D3").Value = "'3.558.56"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'582.28"
$ws.Range("E5").Value = "  +2.17%  "
$ws.Range("D6").Value = "'186.06"
$ws.Range("E6").Value = "  +1.74%  "
$ws.Range("E7").Value = "  +2.72%  "
$ws.Range("D8").Value = "'3.546.29"
$ws.Range("E8").Value = "  +1.96%  "
$ws.Range("E9").Value = "  -0.08%  "
$ws.Range("D10").Value = "'0.220"
$ws.Range("E10").Value = "  +19.96%  "
$ws.Range("E11").Value = "  +1.92%  "
$ws.Range("D12").Value = "'54.43"
$ws.Range("E12").Value = "  +1.57%  "
$ws.Range("D13").Value = "'0.0000316"
$ws.Range("E13").Value = "  +5.77%  "
$ws.Range("E14").Value = "  +0.78%  "
$ws.Range("D15").Value = "'4.125.22"
$ws.Range("E15").Value = "  +1.87%  "
$ws.Range("D16").Value = "'70.820.37"
$ws.Range("E16").Value = "  +3.07%  "
$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").Value = "'19.18"
$ws.Range("E17").Value = "  +0.07%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "'3.535.81"
$ws.Range("E18").Value = "  +1.74%  "
$ws.Range("D19").Value = "'12.49"
$ws.Range("E19").Value = "  +1.65%  "
$ws.Range("D20").Value = "'570.22"
$ws.Range("E20").Value = "  +5.97%  "
$ws.Range("E21").Value = "  +0.81%  "
$ws.Range("E22").Value = "  -0.75%  "
$ws.Range("D23").Value = "'17.65"
$ws.Range("E23").Value = "  -9.04%  "
$ws.Range("D24").Value = "'4.55"
$ws.Range("E24").Value = "  +3.70%  "
$ws.Range("E25").Value = "  -0.99%  "
$ws.Range("D26").Value = "'94.76"
$ws.Range("E26").Value = "  +0.70%  "
$ws.Range("E27").Value = "  +4.73%  "
$ws.Range("E28").Value = "  +1.95%  "
$ws.Range("E29").Value = "  +1.66%  "
$ws.Range("D30").Value = "'32.43"
$ws.Range("E30").Value = "  +3.40%  "
$ws.Range("E31").Value = "  +0.70%  "
$ws.Range("D32").Value = "'12.28"
$ws.Range("E32").Value = "  -1.64%  "
$ws.Range("D33").Value = "'0.116"
$ws.Range("E33").Value = "  +2.60%  "
$ws.Range("D34").Value = "'63.02"
$ws.Range("E34").Value = "  -1.76%  "
$ws.Range("D35").Value = "'3.42"
$ws.Range("E35").Value = "  +12.98%  "
$ws.Range("D36").Value = "'547.32"
$ws.Range("E36").Value = "  -4.21%  "
$ws.Range("D37").Value = "'0.414"
$ws.Range("E37").Value = "  +4.86%  "
$ws.Range("B38").Value = "dogwifhat"
$ws.Range("C38").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D38").Value = "'3.41"
$ws.Range("E38").Value = "  +9.86%  "
$ws.Range("B39").Value = "InjectiveProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D39").Value = "'37.88"
$ws.Range("E39").Value = "  +0.41%  "
$ws.Range("B40").Value = "Dai"
$ws.Range("C40").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D40").Value = "'1.00"
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("B41").Value = "PEPE"
$ws.Range("C41").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D41").Value = "'0.0₃0799"
$ws.Range("E41").Value = "  +5.10%  "
$ws.Range("D42").Value = "'3.592.51"
$ws.Range("E42").Value = "  +11.99%  "
$ws.Range("E43").Value = "  +3.37%  "
$ws.Range("D44").Value = "'3.43"
$ws.Range("E44").Value = "  +3.82%  "
$ws.Range("D45").Value = "'0.0465"
$ws.Range("E45").Value = "  +6.71%  "
$ws.Range("D46").Value = "'3.47"
$ws.Range("E46").Value = "  +1.02%  "
$ws.Range("E47").Value = "  -1.50%  "
$ws.Range("E48").Value = "  +3.50%  "
$ws.Range("E49").Value = "  +3.08%  "
$ws.Range("E50").Value = "  +14.38%  "
$ws.Range("B51").Value = "FirstDigitalUSD"
$ws.Range("C51").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D51").Value = "'0.999"
$ws.Range("E51").Value = "  +0.10%  "
